$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date column C (13-01-2023), matching header formatting of column B
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Rewrite rows 2-52: new sort order (funds alphabetically, then avg, then total) plus new column C values
$ws.Cells.Item(2, 1).Value = "1810 Renta variable"
$ws.Cells.Item(2, 2).Value = 226089.51
$ws.Cells.Item(2, 3).Value = 345962.63
$ws.Cells.Item(3, 1).Value = "1822 Raices Valores Negociables"
$ws.Cells.Item(3, 2).Value = 879586.27
$ws.Cells.Item(3, 3).Value = 879676.12
$ws.Cells.Item(4, 1).Value = "Adcap IOL Acciones Argentina"
$ws.Cells.Item(4, 2).Value = 105473.76
$ws.Cells.Item(4, 3).Value = 105611.31
$ws.Cells.Item(5, 1).Value = "Allaria Acciones"
$ws.Cells.Item(5, 2).Value = 235786.46
$ws.Cells.Item(5, 3).Value = 235983.86
$ws.Cells.Item(6, 1).Value = "Alpha Acciones"
$ws.Cells.Item(6, 2).Value = 286942.29
$ws.Cells.Item(6, 3).Value = 287401.77
$ws.Cells.Item(7, 1).Value = "Alpha Mega"
$ws.Cells.Item(7, 2).Value = 1001323.34
$ws.Cells.Item(7, 3).Value = 1001834.36
$ws.Cells.Item(8, 1).Value = "Alpha Mercosur"
$ws.Cells.Item(8, 2).Value = 598915.4399999999
$ws.Cells.Item(8, 3).Value = 598348.37
$ws.Cells.Item(9, 1).Value = "Alpha Recursos Naturales"
$ws.Cells.Item(9, 2).Value = 756208.8199999999
$ws.Cells.Item(9, 3).Value = 878817.86
$ws.Cells.Item(10, 1).Value = "Alpha planeam equil"
$ws.Cells.Item(10, 2).Value = 21616.12
$ws.Cells.Item(10, 3).Value = 11558.68
$ws.Cells.Item(11, 1).Value = "Alpha renta balan global"
$ws.Cells.Item(11, 2).Value = 944897.53
$ws.Cells.Item(11, 3).Value = 941568.24
$ws.Cells.Item(12, 1).Value = "Argenfunds"
$ws.Cells.Item(12, 2).Value = 35514.84
$ws.Cells.Item(12, 3).Value = 35552.2
$ws.Cells.Item(13, 1).Value = "Arpenta acciones"
$ws.Cells.Item(13, 2).Value = 8193.27
$ws.Cells.Item(13, 3).Value = 8200.26
$ws.Cells.Item(14, 1).Value = "Arpenta ex Mercosur"
$ws.Cells.Item(14, 2).Value = 13408.81
$ws.Cells.Item(14, 3).Value = 13415.14
$ws.Cells.Item(15, 1).Value = "Balanz"
$ws.Cells.Item(15, 2).Value = 699372.66
$ws.Cells.Item(15, 3).Value = 700664.96
$ws.Cells.Item(16, 1).Value = "CMA acciones"
$ws.Cells.Item(16, 2).Value = 196743.44
$ws.Cells.Item(16, 3).Value = 196842.63
$ws.Cells.Item(17, 1).Value = "Compass Crecimiento"
$ws.Cells.Item(17, 2).Value = 2044193.98
$ws.Cells.Item(17, 3).Value = 2035251.4
$ws.Cells.Item(18, 1).Value = "Consultatio Acciones Argentina"
$ws.Cells.Item(18, 2).Value = 1806775.34
$ws.Cells.Item(18, 3).Value = 1681454.57
$ws.Cells.Item(19, 1).Value = "Consultatio Renta Variable"
$ws.Cells.Item(19, 2).Value = 759610.45
$ws.Cells.Item(19, 3).Value = 759475.03
$ws.Cells.Item(20, 1).Value = "Delta Acciones"
$ws.Cells.Item(20, 2).Value = 57514.45
$ws.Cells.Item(20, 3).Value = 57526.56
$ws.Cells.Item(21, 1).Value = "Delta Internacional"
$ws.Cells.Item(21, 2).Value = 15998.35
$ws.Cells.Item(21, 3).Value = 15993.69
$ws.Cells.Item(22, 1).Value = "Delta Latinoamerica"
$ws.Cells.Item(22, 2).Value = 29750.4
$ws.Cells.Item(22, 3).Value = 29730.25
$ws.Cells.Item(23, 1).Value = "Delta Select"
$ws.Cells.Item(23, 2).Value = 477784.31
$ws.Cells.Item(23, 3).Value = 477734.94
$ws.Cells.Item(24, 1).Value = "FBA Acciones Argentinas"
$ws.Cells.Item(24, 2).Value = 611641.45
$ws.Cells.Item(24, 3).Value = 622384.3
$ws.Cells.Item(25, 1).Value = "FBA Calificado"
$ws.Cells.Item(25, 2).Value = 601642.11
$ws.Cells.Item(25, 3).Value = 613750.14
$ws.Cells.Item(26, 1).Value = "Fima Acciones"
$ws.Cells.Item(26, 2).Value = 613831.85
$ws.Cells.Item(26, 3).Value = 864687.73
$ws.Cells.Item(27, 1).Value = "Fima PB Acciones"
$ws.Cells.Item(27, 2).Value = 215635.15
$ws.Cells.Item(27, 3).Value = 458092.97
$ws.Cells.Item(28, 1).Value = "Gainvest Renta Variable"
$ws.Cells.Item(28, 2).Value = 95932.97
$ws.Cells.Item(28, 3).Value = 95762.06
$ws.Cells.Item(29, 1).Value = "Galileo Acciones"
$ws.Cells.Item(29, 2).Value = 1943145.83
$ws.Cells.Item(29, 3).Value = 1942133.3
$ws.Cells.Item(30, 1).Value = "Goal Acciones Argentinas"
$ws.Cells.Item(30, 2).Value = 69522.91
$ws.Cells.Item(30, 3).Value = 69504.27
$ws.Cells.Item(31, 1).Value = "Goal acciones plus"
$ws.Cells.Item(31, 2).Value = 24510.69
$ws.Cells.Item(31, 3).Value = 24517.51
$ws.Cells.Item(32, 1).Value = "HF Acciones Argentinas"
$ws.Cells.Item(32, 2).Value = 330422.7
$ws.Cells.Item(32, 3).Value = 330584.64
$ws.Cells.Item(33, 1).Value = "HF Acciones Lideres"
$ws.Cells.Item(33, 2).Value = 584953.8
$ws.Cells.Item(33, 3).Value = 563917.42
$ws.Cells.Item(34, 1).Value = "IAM Renta Variable"
$ws.Cells.Item(34, 2).Value = 114976.16
$ws.Cells.Item(34, 3).Value = 120265.95
$ws.Cells.Item(35, 1).Value = "IEB Value"
$ws.Cells.Item(35, 2).Value = 22658.35
$ws.Cells.Item(35, 3).Value = 22642.85
$ws.Cells.Item(36, 1).Value = "Lombardi"
$ws.Cells.Item(36, 2).Value = 105680.69
$ws.Cells.Item(36, 3).Value = 105793.94
$ws.Cells.Item(37, 1).Value = "MAF"
$ws.Cells.Item(37, 2).Value = 248774.94
$ws.Cells.Item(37, 3).Value = 248946.29
$ws.Cells.Item(38, 1).Value = "Megainver"
$ws.Cells.Item(38, 2).Value = 99369.00999999999
$ws.Cells.Item(38, 3).Value = 99476.06
$ws.Cells.Item(39, 1).Value = "Pellegrini Acciones"
$ws.Cells.Item(39, 2).Value = 220924.8
$ws.Cells.Item(39, 3).Value = 221233.55
$ws.Cells.Item(40, 1).Value = "Pionero Acciones"
$ws.Cells.Item(40, 2).Value = 399231.85
$ws.Cells.Item(40, 3).Value = 398675.83
$ws.Cells.Item(41, 1).Value = "Premier Renta Variable"
$ws.Cells.Item(41, 2).Value = 142507.79
$ws.Cells.Item(41, 3).Value = 142655.19
$ws.Cells.Item(42, 1).Value = "Quinquela Acciones"
$ws.Cells.Item(42, 2).Value = 320189.3
$ws.Cells.Item(42, 3).Value = 319994.27
$ws.Cells.Item(43, 1).Value = "Rofex 20 Renta Variable"
$ws.Cells.Item(43, 2).Value = 229972.21
$ws.Cells.Item(43, 3).Value = 229555.44
$ws.Cells.Item(44, 1).Value = "SBS Acciones Argentina"
$ws.Cells.Item(44, 2).Value = 1061365.07
$ws.Cells.Item(44, 3).Value = 1061127.24
$ws.Cells.Item(45, 1).Value = "Schroeder RV"
$ws.Cells.Item(45, 2).Value = 2507588.44
$ws.Cells.Item(45, 3).Value = 2506691.17
$ws.Cells.Item(46, 1).Value = "Supefondo RV"
$ws.Cells.Item(46, 2).Value = 1507953.08
$ws.Cells.Item(46, 3).Value = 1513462.28
$ws.Cells.Item(47, 1).Value = "Superfondo "
$ws.Cells.Item(47, 2).Value = 154579.37
$ws.Cells.Item(47, 3).Value = 154269.06
$ws.Cells.Item(48, 1).Value = "Supergestion"
$ws.Cells.Item(48, 2).Value = 319272.67
$ws.Cells.Item(48, 3).Value = 320424.99
$ws.Cells.Item(49, 1).Value = "Toronto Trust Multimercado"
$ws.Cells.Item(49, 2).Value = 114607.64
$ws.Cells.Item(49, 3).Value = 114605.83
$ws.Cells.Item(50, 1).Value = "Toronto trust Argy"
$ws.Cells.Item(50, 2).Value = 189731.39
$ws.Cells.Item(50, 3).Value = 189512.06
$ws.Cells.Item(51, 1).Value = "avg"
$ws.Cells.Item(51, 2).Value = 490863.72
$ws.Cells.Item(51, 3).Value = 503127.98
$ws.Cells.Item(52, 1).Value = "total"
$ws.Cells.Item(52, 2).Value = 24052322.06
$ws.Cells.Item(52, 3).Value = 24653271.17
